# DMo and Com script update
# final updates for CEDARS formatted 8760 output
#
# Renames the "Cap-Tons" Normunit label to "Cap-Ton" everywhere it is used
# (DMo + MFm sheets), and on the Com sheet moves the Cap-Ton conversion
# formulas from column D into column F (where the rest of the "Value"
# column lives) and stamps a Msr label ("SWXX111-00 Example_SEER_AC") in
# the newly-used column G.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DMo: A6 Normunit label Cap-Tons -> Cap-Ton
# ---------------------------------------------------------------------
$wsDMo = $wb.Worksheets.Item("DMo")
$wsDMo.Range("A6").Value = "Cap-Ton"

# ---------------------------------------------------------------------
# MFm: every Normunit cell in column A that reads Cap-Tons -> Cap-Ton
# (rows 5-36 and 38-53; row 37 is the unrelated "kWhreduced" row)
# ---------------------------------------------------------------------
$wsMFm = $wb.Worksheets.Item("MFm")
$mfmRows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53)
foreach ($r in $mfmRows) {
    $wsMFm.Cells.Item($r, 1).Value = "Cap-Ton"
}

# ---------------------------------------------------------------------
# Com: move the Cap-Ton Value formulas from column D to column F, and
# record the measure name in the newly used column G (shares column 7
# with a wider custom width to fit the text).
# ---------------------------------------------------------------------
$wsCom = $wb.Worksheets.Item("Com")

$wsCom.Range("F29").Formula = "=1000512.63*0.000284345"
$wsCom.Range("D29").Value = ""
$wsCom.Range("G29").Value = "SWXX111-00 Example_SEER_AC"

$wsCom.Range("F30").Formula = "=(1028872.44+727611.22+733287.34+56277.3+122055.18+130350.16)*0.000284345"
$wsCom.Range("D30").Value = ""
$wsCom.Range("G30").Value = "SWXX111-00 Example_SEER_AC"

$wsCom.Range("F31").Formula = "=(78241.71+53614.66+65053.12+72956.14+81901.46+70644.79+33603.59)*0.000284345"
$wsCom.Range("D31").Value = ""
$wsCom.Range("G31").Value = "SWXX111-00 Example_SEER_AC"

# widen column G on Com to fit the new measure-name text
$wsCom.Columns.Item(7).ColumnWidth = 28.77734375

# ---------------------------------------------------------------------
# Window / selection state left behind by the author while making the
# edits above (cosmetic, but reproduced for fidelity): DMo, MFm and Com
# each have a new active cell, and SFm ends up the selected/active tab.
# ---------------------------------------------------------------------
$wsSFm = $wb.Worksheets.Item("SFm")

$wsDMo.Range("A6").Select()
$wsMFm.Range("J19").Select()
$wsCom.Range("G41").Select()
$wsSFm.Range("D84").Select()
$wsSFm.Activate()
